# Commit: "Se agrega la operacion Modulo, y se resuelven problemas de
# identacion que tenia el codigo, ya que no arojaba los resultados"
#
# Data-wise this appends a new batch of "El Kevin" submission rows to the
# usage log sheet (rows 14-46) and nudges the last existing timestamp
# (B13) to its corrected value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-use the same date/time number format already applied to the existing
# timestamp column (B2:B13), so the new cells stay visually consistent.
$dateFmt = $ws.Range("B13").NumberFormat

# Correct the slightly-off timestamp that was already in the sheet.
$ws.Cells.Item(13, 2).Value = 45735.92259413195

# New rows of "El Kevin" activity to append below the existing data.
$newTimestamps = @(
    45735.9544002662,
    45735.95449295139,
    45735.9545329051,
    45735.95455835648,
    45735.95455964121,
    45735.95650885416,
    45735.95661969907,
    45735.95668149306,
    45735.95671054398,
    45735.95672408565,
    45735.95673371528,
    45735.95674752315,
    45735.95675637732,
    45735.95735918982,
    45735.95736503472,
    45735.95737634259,
    45735.9573828588,
    45735.95780239583,
    45735.95784375,
    45735.95787296296,
    45735.9578852662,
    45735.95789425926,
    45735.95790657408,
    45735.95791350694,
    45735.95822385417,
    45735.9582303588,
    45735.95823895833,
    45735.95828128472,
    45735.95909799769,
    45735.95920408565,
    45735.95926295139,
    45735.95927793982,
    45735.95928535519
)

$row = 14
foreach ($ts in $newTimestamps) {
    $ws.Cells.Item($row, 1).Value = "El Kevin"

    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = $ts
    $cell.NumberFormat = $dateFmt

    $row = $row + 1
}
